# "Added new actions and results / Matrix fixed 1\5"
#
# The action/result numbering in column B ("Матрица" labels) was out of
# sync with the actual rows of content in column C. This fixes the
# R2/R3/R4 sub-action numbering so every action is sequentially labelled
# (A2.1-A2.5, A3.1-A3.5, A4.1-A4.5) and fills in the labels that were
# previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- R2 block (rows 11-16): Action labels A2.1 .. A2.5 ---
$ws.Range("B12").Value = "A2.1"
$ws.Range("B13").Value = "A2.2"
$ws.Range("B14").Value = "A2.3"
$ws.Range("B15").Value = "A2.4"
$ws.Range("B16").Value = "A2.5"

# --- R3 block (rows 17-22): Action labels A3.1 .. A3.5 ---
$ws.Range("B18").Value = "A3.1"
$ws.Range("B19").Value = "A3.2"
$ws.Range("B20").Value = "A3.3"
$ws.Range("B21").Value = "A3.4"
$ws.Range("B22").Value = "A3.5"

# --- R4 block (rows 23-28): Action labels A4.1 .. A4.5 ---
$ws.Range("B24").Value = "A4.1"
$ws.Range("B25").Value = "A4.2"
$ws.Range("B26").Value = "A4.3"
$ws.Range("B27").Value = "A4.4"
$ws.Range("B28").Value = "A4.5"

# Refresh the selection to match the author's last editing position.
$ws.Range("C23").Select()
